# Append a new "2021年" row (row 11) to the single data sheet, extending the
# yearly time-series table (rows 2-10 hold 2012年..2020年) by one more year.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A10's formatting (bold, centered, thin-bordered "year label" style) onto
# A11 before writing its value, so the new label matches A2:A10 visually.
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A11").Value = "2021年"
$ws.Range("B11").Value = 7043.98
$ws.Range("C11").Value = 1584.88
$ws.Range("D11").Value = 988.23

# E11 must hold an explicit empty-text value (same as E6:E10), not a blank
# cell. Writing "" is treated as "clear the cell", so instead write a lone
# apostrophe (forces an empty text entry) and then strip the transient
# "quote prefix" number-format it introduces, leaving a plain empty string.
$ws.Range("E11").Value = "'"
$ws.Range("E11").ClearFormats()

$ws.Range("F11").Value = 1621.73
$ws.Range("G11").Value = 8942.030000000001
$ws.Range("H11").Value = 683.12
$ws.Range("I11").Value = 4103.39
$ws.Range("J11").Value = 591.23
$ws.Range("K11").Value = 212832.81
$ws.Range("L11").Value = 139.23
$ws.Range("M11").Value = 403
$ws.Range("N11").Value = 1344.33
$ws.Range("O11").Value = 293.16
$ws.Range("P11").Value = 7248.59
$ws.Range("Q11").Value = 1095.06
$ws.Range("R11").Value = 128.37
$ws.Range("S11").Value = 1026.93
$ws.Range("T11").Value = 5866.8
$ws.Range("U11").Value = 21142.16
$ws.Range("V11").Value = 8345.049999999999
$ws.Range("W11").Value = 22052.83
$ws.Range("X11").Value = 2431.06
$ws.Range("Y11").Value = 28367.92
$ws.Range("Z11").Value = 6131.77
$ws.Range("AA11").Value = 54.46
$ws.Range("AB11").Value = 8080.53
$ws.Range("AC11").Value = 2942.46
$ws.Range("AD11").Value = 497.09
$ws.Range("AE11").Value = 257.83
$ws.Range("AF11").Value = 13542.46
$ws.Range("AG11").Value = 8662.98
$ws.Range("AH11").Value = 733.48
$ws.Range("AI11").Value = 5408.87
$ws.Range("AJ11").Value = 879.15
$ws.Range("AK11").Value = 3536.25
$ws.Range("AL11").Value = 15300.64
$ws.Range("AM11").Value = 7257.35
$ws.Range("AN11").Value = 1247.99
$ws.Range("AO11").Value = 975.52
$ws.Range("AP11").Value = 10007.66
$ws.Range("AQ11").Value = 1870.31
